$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 509; this pushes the old rows 509-514 down to 513-518
$ws.Rows.Item(509).Resize(4).Insert()

# Common/constant values for this block of records
$mercadoId = 5
$mercado = "Macroferia Regional de Talca"
$region = "Maule"
$codreg = 7
$tipo = "Fruta"
$productoId = 100104
$producto = "Frutos de pepita"
$categoriaId = 100104005
$categoria = "Pera"
$unidad = "`$/bandeja 18 kilos granel"
$kgUnidad = 18

function Set-Row {
    param(
        [int]$r,
        [string]$variedad,
        [string]$calidad,
        [double]$fecha,
        [double]$volumen,
        [double]$precioMin,
        [double]$precioMax,
        [double]$precioProm,
        [string]$origen,
        [double]$precioKg
    )

    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $tipo
    $ws.Cells.Item($r, 7).Value = $productoId
    $ws.Cells.Item($r, 8).Value = $producto
    $ws.Cells.Item($r, 9).Value = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $calidad
    $ws.Cells.Item($r, 13).Value = $volumen
    $ws.Cells.Item($r, 14).Value = $precioMin
    $ws.Cells.Item($r, 15).Value = $precioMax
    $ws.Cells.Item($r, 16).Value = $precioProm
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $precioKg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}

Set-Row 509 "Forelle"              "Primera"  44656 180 8000  8000  8000  "Provincia de Curicó" 444
Set-Row 510 "Packham's Triumph"    "Especial" 44656 230 10000 10000 10000 "Provincia de Curicó" 556
Set-Row 511 "Packham's Triumph"    "Especial" 44656 260 10000 10000 10000 "Provincia de Linares" 556
Set-Row 512 "Packham's Triumph"    "Primera"  44656 200 8000  8000  8000  "Provincia de Curicó" 444
